# v2.6 Added decoupled suspension, four-wheel steering, scripts to generate GGV diagram
#
# Adds a new "Semi_Truck_Scalable" sheet to the 3-axle body library by
# duplicating the existing "Truck_Amandla_3Axle" sheet (same layout/
# styles/formulas), renaming it, relabelling its header cell, and making
# it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# The new sheet is a duplicate of the Truck sheet (same cols/rows/styles/
# formulas) - mirrors how this was authored in Excel (copy sheet -> rename).
$srcWs = $wb.Worksheets.Item("Truck_Amandla_3Axle")
$srcWs.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "Semi_Truck_Scalable"

# Row 3 / column H carries the sheet's own name as a label - update it to
# match the new sheet.
$newWs.Range("H3").Value = "Semi_Truck_Scalable"

# Leave the new sheet selected/active, with the last-used cell selected.
$newWs.Activate()
$newWs.Range("J17").Select()

# The Truck sheet's own last selection moved too.
$srcWs.Activate()
$srcWs.Range("D24").Select()

# Re-activate the new sheet so it ends up the active/visible tab.
$newWs.Activate()
